# "some ui art for magmaCooler scene"
# Adds new localisation rows to the "en" language sheet:
#   - a block of 5 new key/value rows (select_method, intrusive, extrusive,
#     back, proceed) inserted right before the existing "victory" row
#   - a block of 2 new key/value rows (cooling, stop) inserted right before
#     the existing "grainSize_LargeVariant" row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first block: 5 rows above current row 15 ("victory") ---------
$ws.Range("A15:A19").EntireRow.Insert()

$ws.Range("A15").Value() = "select_method"
$ws.Range("B15").Value() = "Select Method"

$ws.Range("A16").Value() = "intrusive"
$ws.Range("B16").Value() = "Intrusive"

$ws.Range("A17").Value() = "extrusive"
$ws.Range("B17").Value() = "Extrusive"

$ws.Range("A18").Value() = "back"
$ws.Range("B18").Value() = "BACK"

$ws.Range("A19").Value() = "proceed"
$ws.Range("B19").Value() = "PROCEED"

# --- Insert second block: 2 rows above current row 27
#     ("grainSize_LargeVariant", after the shift from the first insert) ----
$ws.Range("A27:A28").EntireRow.Insert()

$ws.Range("A27").Value() = "cooling"
$ws.Range("A28").Value() = "stop"
$ws.Range("B28").Value() = "STOP"
$ws.Range("B27").Value() = "Cooling…"

# --- Cursor / selection bookkeeping (cosmetic, matches author's session) --
$ws.Range("B27").Select()
